$d = $word.ActiveDocument

$LQ = [char]0x201C
$RQ = [char]0x201D

function Add-BulletParagraph {
    param(
        [int]$Level,
        [string]$Text,
        [bool]$MarkFormatted = $false,
        [string]$BookmarkName = ""
    )

    $end = $d.Content.End
    $insertionRange = $d.Range($end, $end)
    $insertionRange.InsertParagraphAfter()

    $newPara = $d.Paragraphs.Last
    $newPara.Range.ListFormat.ListLevelNumber = $Level

    $textRange = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)
    $textRange.Text = $Text

    if ($BookmarkName) {
        $startPoint = $d.Range($newPara.Range.Start, $newPara.Range.Start)
        $d.Bookmarks.Add($BookmarkName, $startPoint)
    }

    if ($MarkFormatted) {
        $newPara.Range.Font.Size = 12
        $newPara.Range.Font.SizeBi = 12
    }

    return $newPara
}

# --- Gard et al., 2018 ---
Add-BulletParagraph 2 "Gard et al., 2018" | Out-Null
Add-BulletParagraph 3 "Issues:" | Out-Null
Add-BulletParagraph 4 "Reason for operationalization given in operationalization quote not given in reason quote" | Out-Null
Add-BulletParagraph 3 "Resolution:" | Out-Null
Add-BulletParagraph 4 "Copied operationalization quote into reason cell" $false "__DdeLink__9209_101083012" | Out-Null

# --- Bayissa et al., 2017 ---
Add-BulletParagraph 2 "Bayissa et al., 2017" | Out-Null
Add-BulletParagraph 3 "Issues:" | Out-Null
$t8 = "RA did not update to " + $LQ + "Did not study SES" + $RQ + " per my comments"
Add-BulletParagraph 4 $t8 | Out-Null
Add-BulletParagraph 3 "Resolution:" | Out-Null
$t10 = "Deleted all codes and made " + $LQ + "sesIVorDV" + $RQ + " " + $LQ + "Did not study SES" + $RQ
Add-BulletParagraph 4 $t10 | Out-Null

# --- Carlson et al., 2014 ---
Add-BulletParagraph 2 "Carlson et al., 2014" | Out-Null
Add-BulletParagraph 3 "Issues:" | Out-Null
$t13 = "Reason for operationalization in operationalization quote was not in reason quote (which was " + $LQ + "None given" + $RQ + ")"
Add-BulletParagraph 4 $t13 | Out-Null
Add-BulletParagraph 3 "Resolution:" $true | Out-Null
Add-BulletParagraph 4 "Copied operationalization quote into reason cell" $true | Out-Null

Write-Output ("Final paragraph count=" + $d.Paragraphs.Count)
